$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename existing "TestN" placeholder asset names to real ones ---
$ws.Range("D2").Value = "Character01"
$ws.Range("F2").Value = "Character01_Image"

$ws.Range("D3").Value = "Gun01"
$ws.Range("F3").Value = "Gun01_Image"

$ws.Range("D4").Value = "UpperArmor01"
$ws.Range("F4").Value = "UpperArmor01_Image"

# --- Update Id column values to the new large numeric ids ---
$ws.Range("A2").Value = 10000001
$ws.Range("A3").Value = 10000002
$ws.Range("A4").Value = 10000003

# --- Widen column F to fit the longer asset names ---
$ws.Range("F1").EntireColumn.ColumnWidth = 24.5

# --- Add new row 5: Item / HeathPotion01 (consumable) ---
# Copy the formatting of row 4 first so the new row matches existing styling.
$ws.Range("A4:E4").Copy()
$ws.Range("A5:E5").PasteSpecial(-4122)
$ws.Range("E4").Copy()
$ws.Range("F5").PasteSpecial(-4122)

$ws.Range("A5").Value = 10000004
$ws.Range("B5").Value = "Item"
$ws.Range("C5").Value = "Assets/Prefabs/Item"
$ws.Range("D5").Value = "HeathPotion01"
$ws.Range("E5").Value = "Assets/Images/"
$ws.Range("F5").Value = "HeathPotion01_Image"
